# Add two new columns, I ("I0") and J ("IF"), to the sheet.
#   - I0: constant 1 for every data row.
#   - IF: same value as column H for every data row.
# Header cells (row 1) get the same (bold / centered / bordered) style as
# the existing headers in B1:H1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -------------------------------------------------------------
# Copy the header style from H1 (formats only) onto I1 and J1, then set
# their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows -------------------------------------------------------------
$lastRow = 22
for ($r = 2; $r -le $lastRow; $r++) {
    $hValue = $ws.Cells.Item($r, 8).Value2   # column H
    $ws.Cells.Item($r, 9).Value = 1           # column I (I0)
    $ws.Cells.Item($r, 10).Value = $hValue    # column J (IF)
}
